# "Generate Report for Archive"
#
# The localization-status report re-sorted the handoff table: the row for
# file 69122473-a34a-4d10-b38d-1fdd26acc32e (previously the 4th data row,
# sheet row 8) now sorts as the 1st of that group (sheet row 5). The three
# rows that used to precede it (86f42771-..., fb659db4-..., 0689cad6-...)
# each shift down by one row. This happens identically on all three
# worksheets (Overview, zh-cn, de-de). Everything else (header row, rows
# 2-4, row 9, styles, hyperlink rIds) is untouched.

$wb = $excel.ActiveWorkbook

foreach ($wsName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($wsName)

    $lastCol = $ws.UsedRange.Columns.Count

    # Snapshot the 4 affected rows (sheet rows 5..8) across every used column.
    $snapshot = @()
    for ($r = 5; $r -le 8; $r++) {
        $rowVals = @()
        for ($c = 1; $c -le $lastCol; $c++) {
            $rowVals += ,$ws.Cells.Item($r, $c).Value2
        }
        $snapshot += ,$rowVals
    }

    # snapshot[0] = old row5 (86f42771-...)
    # snapshot[1] = old row6 (fb659db4-...)
    # snapshot[2] = old row7 (0689cad6-...)
    # snapshot[3] = old row8 (69122473-...)
    #
    # New order: row5 <- old row8, row6 <- old row5, row7 <- old row6, row8 <- old row7
    $newOrder = @($snapshot[3], $snapshot[0], $snapshot[1], $snapshot[2])

    for ($i = 0; $i -lt 4; $i++) {
        $destRow = 5 + $i
        $rowVals = $newOrder[$i]
        for ($c = 1; $c -le $lastCol; $c++) {
            $ws.Cells.Item($destRow, $c).Value2 = $rowVals[$c - 1]
        }
    }
}
